$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "2026-02-03 Tuesday" "2026-02-04 Wednesday"

Replace-Text "318×6=" "958×6="
Replace-Text "986×6=" "975×2="
Replace-Text "316×9=" "830×2="
Replace-Text "922×6=" "976×9="
Replace-Text "337×2=" "162×7="

Replace-Text "539×3=" "781×8="
Replace-Text "528×5=" "697×7="
Replace-Text "616×2=" "881×8="
Replace-Text "557×4=" "847×2="
Replace-Text "625×9=" "677×3="

Replace-Text "155×8=" "864×4="
Replace-Text "287×6=" "490×4="
Replace-Text "427×5=" "988×9="
Replace-Text "207×6=" "934×4="
Replace-Text "125×5=" "996×8="

Replace-Text "532×3=" "361×2="
Replace-Text "309×2=" "770×6="
Replace-Text "867×2=" "326×5="
Replace-Text "858×2=" "808×4="
Replace-Text "944×7=" "111×7="

Replace-Text "518×8=" "348×5="
Replace-Text "522×7=" "879×2="
Replace-Text "421×9=" "880×8="
Replace-Text "316×6=" "548×6="
Replace-Text "991×3=" "676×7="
